$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new quantity figures in columns H and I for a few components
$ws.Range("I7").Value = 50

$ws.Range("H8").Value = 105
$ws.Range("I8").Formula = "=SUM(H8/4)"

$ws.Range("I14").Value = 58

# Hide the now-secondary lookup/reference columns C:G
$ws.Range("C:G").EntireColumn.Hidden = $true

# Update the active selection to reflect where the user ended up working
$ws.Range("K12").Select() | Out-Null
